$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3 was stored as an inline string "9087654321"; convert it to a real number
$ws.Range("E3").Value = 9087654321

# Add new row 4
$ws.Range("A4").Value = "svbhadri1110@gmail.com"
$ws.Range("B4").Value = "passwd123"
$ws.Range("C4").Value = "MALE"
$ws.Range("D4").Value = 21
$ws.Range("E4").Value = 9087654321

# Add new row 5
$ws.Range("A5").Value = "asbhj@fkjas.com"
$ws.Range("B5").Value = "efwfewfe"
$ws.Range("C5").Value = "MALE"
$ws.Range("D5").Value = 20
# E5 stays a text string "1234567890" (not converted to a number)
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1234567890"
